# Update margin figures on the "GPC" worksheet to reflect the revised
# income-statement percentages (rows 15-19 and 30-31, columns D-G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GPC")

# Row 15 - Gross Margin
$ws.Range("D15").Value = 0.3308
$ws.Range("E15").Value = 0.3305
$ws.Range("F15").Value = 0.3333
$ws.Range("G15").Value = 0.3344

# Row 16 - EBIT Margin
$ws.Range("D16").Value = 0.0183
$ws.Range("E16").Value = 0.0158
$ws.Range("F16").Value = 0.0485
$ws.Range("G16").Value = 0.052

# Row 17 - EBT margin
$ws.Range("D17").Value = 0.0154
$ws.Range("E17").Value = 0.0144
$ws.Range("F17").Value = 0.0477
$ws.Range("G17").Value = 0.049

# Row 18 - Net Profit Margin
$ws.Range("E18").Value = 0.0011
$ws.Range("F18").Value = 0.0345
$ws.Range("G18").Value = 0.0354

# Row 19 - Free Cash Flow Margin
$ws.Range("D19").Value = 0.0854
$ws.Range("E19").Value = 0.0783
$ws.Range("F19").Value = 0.0345
$ws.Range("G19").Value = 0.033

# Row 30 - EBITDA Margin
$ws.Range("D30").Value = 0.0707
$ws.Range("E30").Value = 0.0667
$ws.Range("F30").Value = 0.0658
$ws.Range("G30").Value = 0.0667

# Row 31 - Operating Cash Flow Margin
$ws.Range("D31").Value = 0.0964
$ws.Range("E31").Value = 0.092
$ws.Range("F31").Value = 0.0487
$ws.Range("G31").Value = 0.0475
